# Append two new game-history rows to the "game history" sheet
# (sheet2 = "היסטוריית משחקים"), recording two new game results that
# were previously tracked against a delay (against the clock / against
# the computer), matching the commit "Computer pick delay removed /
# Downgraded computers turn pick delay".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("היסטוריית משחקים")

# Row 40: קיריל צ'רננקו, played "נגד הזמן" (against the clock)
$ws.Range("A40").Value = "קיריל צ'רננקו"
$ws.Range("B40").Value = 43992.93096039352
$ws.Range("B40").NumberFormat = "dd-MM-yyyy"
$ws.Range("C40").Value = "נגד הזמן"
$ws.Range("D40").Value = 115

# Row 41: ליאת נתח, played "נגד המחשב" (against the computer)
$ws.Range("A41").Value = "ליאת נתח"
$ws.Range("B41").Value = 43992.93201335648
$ws.Range("B41").NumberFormat = "dd-MM-yyyy"
$ws.Range("C41").Value = "נגד המחשב"
$ws.Range("D41").Value = 70
